# Reworked how the econ data is read and lagged; add pre-2020 / 2020-beyond
# train-test split helper column: a 2008-9 recession dummy variable is
# inserted as a new column Q on Sheet1 (pushing the existing
# "AnnualizedMoM-CPI-Inflation" column from Q to R), and the "Graphs" sheet's
# Chart 15 series (which plotted the old Q column) is re-pointed at R.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Insert a new column at Q (17), shifting the old Q (and everything to its
#    right) one column to the right.
$ws.Columns.Item(17).Insert()

# 2) Header for the freshly inserted column.
$ws.Range("Q1").Value = "2008-9RecessionDummyVar"

# 3) Fill the dummy variable down Q2:Q281 - 1 for the 2008-9 recession
#    months (Jul 2008 - Jan 2009, rows 104-110), 0 everywhere else.
for ($r = 2; $r -le 281; $r++) {
    if ($r -ge 104 -and $r -le 110) {
        $ws.Cells.Item($r, 17).Value = 1
    } else {
        $ws.Cells.Item($r, 17).Value = 0
    }
}

# 4) The chart on the Graphs sheet ("Chart 15") plotted the old Q column;
#    repoint its series formula at the shifted column R.
$graphs = $wb.Worksheets.Item("Graphs")
$chart = $graphs.ChartObjects().Item(15).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Sheet1!`$R`$1,,Sheet1!`$R`$2:`$R`$424,1)"

# 5) Freeze the first column and leave the selection on the new Q2 cell, as
#    in the saved workbook.
$ws.Activate()
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("Q2").Select()
